# Updates the worksheet date and regenerates the 100 addition/subtraction
# problems in the table, matching the "Update master to output generated
# at 4250d90" commit. Each cell's text is unique in the document, so a
# plain literal Find/Replace (no wildcards) on $d.Content is safe and
# unambiguous for every one of the 101 replacements below.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-22 Monday", 2) | Out-Null
$d.Content.Find.Execute("57+11=", $true, $false, $false, $false, $false, $true, 1, $false, "96-91=", 2) | Out-Null
$d.Content.Find.Execute("20+36=", $true, $false, $false, $false, $false, $true, 1, $false, "80-26=", 2) | Out-Null
$d.Content.Find.Execute("92-0=", $true, $false, $false, $false, $false, $true, 1, $false, "65+25=", 2) | Out-Null
$d.Content.Find.Execute("68-24=", $true, $false, $false, $false, $false, $true, 1, $false, "19+42=", 2) | Out-Null
$d.Content.Find.Execute("13+86=", $true, $false, $false, $false, $false, $true, 1, $false, "57+42=", 2) | Out-Null
$d.Content.Find.Execute("6+34=", $true, $false, $false, $false, $false, $true, 1, $false, "17-2=", 2) | Out-Null
$d.Content.Find.Execute("15+57=", $true, $false, $false, $false, $false, $true, 1, $false, "32-22=", 2) | Out-Null
$d.Content.Find.Execute("11+13=", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=", 2) | Out-Null
$d.Content.Find.Execute("64+23=", $true, $false, $false, $false, $false, $true, 1, $false, "48+10=", 2) | Out-Null
$d.Content.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "89+1=", 2) | Out-Null
$d.Content.Find.Execute("52-35=", $true, $false, $false, $false, $false, $true, 1, $false, "82-1=", 2) | Out-Null
$d.Content.Find.Execute("26+8=", $true, $false, $false, $false, $false, $true, 1, $false, "60+0=", 2) | Out-Null
$d.Content.Find.Execute("99-25=", $true, $false, $false, $false, $false, $true, 1, $false, "54-22=", 2) | Out-Null
$d.Content.Find.Execute("7+43=", $true, $false, $false, $false, $false, $true, 1, $false, "89-5=", 2) | Out-Null
$d.Content.Find.Execute("26-22=", $true, $false, $false, $false, $false, $true, 1, $false, "9+61=", 2) | Out-Null
$d.Content.Find.Execute("8+31=", $true, $false, $false, $false, $false, $true, 1, $false, "5+62=", 2) | Out-Null
$d.Content.Find.Execute("36-3=", $true, $false, $false, $false, $false, $true, 1, $false, "89-30=", 2) | Out-Null
$d.Content.Find.Execute("15-3=", $true, $false, $false, $false, $false, $true, 1, $false, "1+65=", 2) | Out-Null
$d.Content.Find.Execute("78+21=", $true, $false, $false, $false, $false, $true, 1, $false, "78-25=", 2) | Out-Null
$d.Content.Find.Execute("0+33=", $true, $false, $false, $false, $false, $true, 1, $false, "20+26=", 2) | Out-Null
$d.Content.Find.Execute("91-33=", $true, $false, $false, $false, $false, $true, 1, $false, "25+2=", 2) | Out-Null
$d.Content.Find.Execute("66-59=", $true, $false, $false, $false, $false, $true, 1, $false, "60+37=", 2) | Out-Null
$d.Content.Find.Execute("85-61=", $true, $false, $false, $false, $false, $true, 1, $false, "32+9=", 2) | Out-Null
$d.Content.Find.Execute("22+4=", $true, $false, $false, $false, $false, $true, 1, $false, "12+9=", 2) | Out-Null
$d.Content.Find.Execute("93-3=", $true, $false, $false, $false, $false, $true, 1, $false, "62-35=", 2) | Out-Null
$d.Content.Find.Execute("81+6=", $true, $false, $false, $false, $false, $true, 1, $false, "99-5=", 2) | Out-Null
$d.Content.Find.Execute("79-71=", $true, $false, $false, $false, $false, $true, 1, $false, "22+63=", 2) | Out-Null
$d.Content.Find.Execute("31-9=", $true, $false, $false, $false, $false, $true, 1, $false, "0+18=", 2) | Out-Null
$d.Content.Find.Execute("48+14=", $true, $false, $false, $false, $false, $true, 1, $false, "61-45=", 2) | Out-Null
$d.Content.Find.Execute("45-8=", $true, $false, $false, $false, $false, $true, 1, $false, "50+0=", 2) | Out-Null
$d.Content.Find.Execute("28+47=", $true, $false, $false, $false, $false, $true, 1, $false, "74-0=", 2) | Out-Null
$d.Content.Find.Execute("6+82=", $true, $false, $false, $false, $false, $true, 1, $false, "62-36=", 2) | Out-Null
$d.Content.Find.Execute("23-23=", $true, $false, $false, $false, $false, $true, 1, $false, "74-29=", 2) | Out-Null
$d.Content.Find.Execute("76-46=", $true, $false, $false, $false, $false, $true, 1, $false, "45+14=", 2) | Out-Null
$d.Content.Find.Execute("38+40=", $true, $false, $false, $false, $false, $true, 1, $false, "90-65=", 2) | Out-Null
$d.Content.Find.Execute("9+53=", $true, $false, $false, $false, $false, $true, 1, $false, "83-73=", 2) | Out-Null
$d.Content.Find.Execute("9+21=", $true, $false, $false, $false, $false, $true, 1, $false, "52-36=", 2) | Out-Null
$d.Content.Find.Execute("80-16=", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=", 2) | Out-Null
$d.Content.Find.Execute("50+5=", $true, $false, $false, $false, $false, $true, 1, $false, "77+21=", 2) | Out-Null
$d.Content.Find.Execute("31+58=", $true, $false, $false, $false, $false, $true, 1, $false, "44-42=", 2) | Out-Null
$d.Content.Find.Execute("63+20=", $true, $false, $false, $false, $false, $true, 1, $false, "5+11=", 2) | Out-Null
$d.Content.Find.Execute("15+16=", $true, $false, $false, $false, $false, $true, 1, $false, "26-23=", 2) | Out-Null
$d.Content.Find.Execute("90-57=", $true, $false, $false, $false, $false, $true, 1, $false, "54+18=", 2) | Out-Null
$d.Content.Find.Execute("50+20=", $true, $false, $false, $false, $false, $true, 1, $false, "32+18=", 2) | Out-Null
$d.Content.Find.Execute("98-21=", $true, $false, $false, $false, $false, $true, 1, $false, "4+84=", 2) | Out-Null
$d.Content.Find.Execute("65+16=", $true, $false, $false, $false, $false, $true, 1, $false, "20+16=", 2) | Out-Null
$d.Content.Find.Execute("91-57=", $true, $false, $false, $false, $false, $true, 1, $false, "1+9=", 2) | Out-Null
$d.Content.Find.Execute("53+23=", $true, $false, $false, $false, $false, $true, 1, $false, "97-56=", 2) | Out-Null
$d.Content.Find.Execute("60-40=", $true, $false, $false, $false, $false, $true, 1, $false, "35+42=", 2) | Out-Null
$d.Content.Find.Execute("32-10=", $true, $false, $false, $false, $false, $true, 1, $false, "90-61=", 2) | Out-Null
$d.Content.Find.Execute("97-6=", $true, $false, $false, $false, $false, $true, 1, $false, "94-44=", 2) | Out-Null
$d.Content.Find.Execute("31-13=", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=", 2) | Out-Null
$d.Content.Find.Execute("65-45=", $true, $false, $false, $false, $false, $true, 1, $false, "3+35=", 2) | Out-Null
$d.Content.Find.Execute("67-51=", $true, $false, $false, $false, $false, $true, 1, $false, "5+78=", 2) | Out-Null
$d.Content.Find.Execute("77-41=", $true, $false, $false, $false, $false, $true, 1, $false, "69+7=", 2) | Out-Null
$d.Content.Find.Execute("73-29=", $true, $false, $false, $false, $false, $true, 1, $false, "75-34=", 2) | Out-Null
$d.Content.Find.Execute("71+4=", $true, $false, $false, $false, $false, $true, 1, $false, "91-79=", 2) | Out-Null
$d.Content.Find.Execute("69+12=", $true, $false, $false, $false, $false, $true, 1, $false, "3+11=", 2) | Out-Null
$d.Content.Find.Execute("72-8=", $true, $false, $false, $false, $false, $true, 1, $false, "3+57=", 2) | Out-Null
$d.Content.Find.Execute("55+39=", $true, $false, $false, $false, $false, $true, 1, $false, "82-60=", 2) | Out-Null
$d.Content.Find.Execute("99-1=", $true, $false, $false, $false, $false, $true, 1, $false, "96-2=", 2) | Out-Null
$d.Content.Find.Execute("27+20=", $true, $false, $false, $false, $false, $true, 1, $false, "52+7=", 2) | Out-Null
$d.Content.Find.Execute("80-32=", $true, $false, $false, $false, $false, $true, 1, $false, "95-91=", 2) | Out-Null
$d.Content.Find.Execute("78-31=", $true, $false, $false, $false, $false, $true, 1, $false, "99-87=", 2) | Out-Null
$d.Content.Find.Execute("76-25=", $true, $false, $false, $false, $false, $true, 1, $false, "59-56=", 2) | Out-Null
$d.Content.Find.Execute("33+51=", $true, $false, $false, $false, $false, $true, 1, $false, "44-39=", 2) | Out-Null
$d.Content.Find.Execute("22-4=", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=", 2) | Out-Null
$d.Content.Find.Execute("72-46=", $true, $false, $false, $false, $false, $true, 1, $false, "83-9=", 2) | Out-Null
$d.Content.Find.Execute("34+64=", $true, $false, $false, $false, $false, $true, 1, $false, "0+97=", 2) | Out-Null
$d.Content.Find.Execute("27-24=", $true, $false, $false, $false, $false, $true, 1, $false, "62-2=", 2) | Out-Null
$d.Content.Find.Execute("80-5=", $true, $false, $false, $false, $false, $true, 1, $false, "48-32=", 2) | Out-Null
$d.Content.Find.Execute("4+32=", $true, $false, $false, $false, $false, $true, 1, $false, "84-55=", 2) | Out-Null
$d.Content.Find.Execute("5+13=", $true, $false, $false, $false, $false, $true, 1, $false, "13+1=", 2) | Out-Null
$d.Content.Find.Execute("38+10=", $true, $false, $false, $false, $false, $true, 1, $false, "87+8=", 2) | Out-Null
$d.Content.Find.Execute("96-87=", $true, $false, $false, $false, $false, $true, 1, $false, "72-55=", 2) | Out-Null
$d.Content.Find.Execute("58-3=", $true, $false, $false, $false, $false, $true, 1, $false, "49+26=", 2) | Out-Null
$d.Content.Find.Execute("21+60=", $true, $false, $false, $false, $false, $true, 1, $false, "57+19=", 2) | Out-Null
$d.Content.Find.Execute("32+23=", $true, $false, $false, $false, $false, $true, 1, $false, "35-23=", 2) | Out-Null
$d.Content.Find.Execute("58-23=", $true, $false, $false, $false, $false, $true, 1, $false, "38+57=", 2) | Out-Null
$d.Content.Find.Execute("58+6=", $true, $false, $false, $false, $false, $true, 1, $false, "26-5=", 2) | Out-Null
$d.Content.Find.Execute("81+9=", $true, $false, $false, $false, $false, $true, 1, $false, "72-58=", 2) | Out-Null
$d.Content.Find.Execute("8+87=", $true, $false, $false, $false, $false, $true, 1, $false, "70-12=", 2) | Out-Null
$d.Content.Find.Execute("14+74=", $true, $false, $false, $false, $false, $true, 1, $false, "90-68=", 2) | Out-Null
$d.Content.Find.Execute("15-10=", $true, $false, $false, $false, $false, $true, 1, $false, "78-59=", 2) | Out-Null
$d.Content.Find.Execute("37-30=", $true, $false, $false, $false, $false, $true, 1, $false, "55-6=", 2) | Out-Null
$d.Content.Find.Execute("35+39=", $true, $false, $false, $false, $false, $true, 1, $false, "80-22=", 2) | Out-Null
$d.Content.Find.Execute("35+30=", $true, $false, $false, $false, $false, $true, 1, $false, "63-15=", 2) | Out-Null
$d.Content.Find.Execute("90-72=", $true, $false, $false, $false, $false, $true, 1, $false, "20+53=", 2) | Out-Null
$d.Content.Find.Execute("60+4=", $true, $false, $false, $false, $false, $true, 1, $false, "15+53=", 2) | Out-Null
$d.Content.Find.Execute("19-5=", $true, $false, $false, $false, $false, $true, 1, $false, "22+49=", 2) | Out-Null
$d.Content.Find.Execute("4+5=", $true, $false, $false, $false, $false, $true, 1, $false, "75+2=", 2) | Out-Null
$d.Content.Find.Execute("68+1=", $true, $false, $false, $false, $false, $true, 1, $false, "29-16=", 2) | Out-Null
$d.Content.Find.Execute("81+16=", $true, $false, $false, $false, $false, $true, 1, $false, "58+25=", 2) | Out-Null
$d.Content.Find.Execute("92-29=", $true, $false, $false, $false, $false, $true, 1, $false, "58-15=", 2) | Out-Null
$d.Content.Find.Execute("8+8=", $true, $false, $false, $false, $false, $true, 1, $false, "54+42=", 2) | Out-Null
$d.Content.Find.Execute("9+36=", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=", 2) | Out-Null
$d.Content.Find.Execute("84-23=", $true, $false, $false, $false, $false, $true, 1, $false, "88-10=", 2) | Out-Null
$d.Content.Find.Execute("71+24=", $true, $false, $false, $false, $false, $true, 1, $false, "91-73=", 2) | Out-Null
$d.Content.Find.Execute("85-72=", $true, $false, $false, $false, $false, $true, 1, $false, "37+49=", 2) | Out-Null
$d.Content.Find.Execute("53-49=", $true, $false, $false, $false, $false, $true, 1, $false, "49-40=", 2) | Out-Null
